$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7304773333333333
$ws.Range("H2").Value = 2.191432
$ws.Range("I2").Value = 0.03163269997405359
$ws.Range("J2").Value = 0.03163269997405359
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 80.82199029261865
$ws.Range("R2").Value = 727.3979126335679
$ws.Range("S2").Value = 0.01732339167989822
$ws.Range("T2").Value = 0.01732339167989822

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7304773333333333
$ws.Range("H3").Value = 2.191432
$ws.Range("I3").Value = 0.03163269997405359
$ws.Range("J3").Value = 0.03163269997405359
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 46.50410697792444
$ws.Range("R3").Value = 418.53696280132
$ws.Range("S3").Value = 0.009967693903425785
$ws.Range("T3").Value = 0.009967693903425785

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7304773333333333
$ws.Range("H4").Value = 2.191432
$ws.Range("I4").Value = 0.03163269997405359
$ws.Range("J4").Value = 0.03163269997405359
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.188324
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 20.25572835999644
$ws.Range("R4").Value = 182.301555239968
$ws.Range("S4").Value = 0.004341614390729593
$ws.Range("T4").Value = 0.004341614390729593

# Row 5
$ws.Range("I5").Value = 0.4074771110502447
$ws.Range("J5").Value = 0.4074771110502448
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 1041.10971054574
$ws.Range("R5").Value = 9369.987394911659
$ws.Range("S5").Value = 0.2231515362617403
$ws.Range("T5").Value = 0.2231515362617403

# Row 6
$ws.Range("I6").Value = 0.4074771110502447
$ws.Range("J6").Value = 0.4074771110502448
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("S6").Value = 0.1283990022645101
$ws.Range("T6").Value = 0.1283990022645101

# Row 7
$ws.Range("I7").Value = 0.4074771110502447
$ws.Range("J7").Value = 0.4074771110502448
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.188324
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 260.9244762894066
$ws.Range("R7").Value = 2348.32028660466
$ws.Range("S7").Value = 0.05592657252399437
$ws.Range("T7").Value = 0.05592657252399438

# Row 8
$ws.Range("G8").Value = 12.95234266666667
$ws.Range("H8").Value = 38.857028
$ws.Range("I8").Value = 0.5608901889757016
$ws.Range("J8").Value = 0.5608901889757018
$ws.Range("M8").Value = 110.642708
$ws.Range("N8").Value = 331.928124
$ws.Range("O8").Value = 0.5476418925386564
$ws.Range("P8").Value = 0.5476418925386564
$ws.Range("Q8").Value = 1433.082267583941
$ws.Range("R8").Value = 12897.74040825547
$ws.Range("S8").Value = 0.3071669645970179
$ws.Range("T8").Value = 0.3071669645970179

# Row 9
$ws.Range("G9").Value = 12.95234266666667
$ws.Range("H9").Value = 38.857028
$ws.Range("I9").Value = 0.5608901889757016
$ws.Range("J9").Value = 0.5608901889757018
$ws.Range("O9").Value = 0.3151072754333865
$ws.Range("P9").Value = 0.3151072754333865
$ws.Range("Q9").Value = 824.5801772339756
$ws.Range("R9").Value = 7421.221595105781
$ws.Range("S9").Value = 0.1767405792654506
$ws.Range("T9").Value = 0.1767405792654507

# Row 10
$ws.Range("G10").Value = 12.95234266666667
$ws.Range("H10").Value = 38.857028
$ws.Range("I10").Value = 0.5608901889757016
$ws.Range("J10").Value = 0.5608901889757018
$ws.Range("M10").Value = 27.72944133333333
$ws.Range("N10").Value = 83.188324
$ws.Range("O10").Value = 0.1372508320279571
$ws.Range("P10").Value = 0.1372508320279571
$ws.Range("Q10").Value = 359.1612261045635
$ws.Range("R10").Value = 3232.451034941072
$ws.Range("S10").Value = 0.07698264511323316
$ws.Range("T10").Value = 0.07698264511323318
